$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet originally listed 6 contract rows (rows 2-7) plus a totals row (row 8).
# The update reduces this to 2 contract rows (rows 2-3) plus the totals row (now row 4):
#   - Row 2 keeps contract "002/DR002/AV1" but renames it to "002/DR002" and
#     updates its monetary figures (loyer/avance/taxes/caution/net).
#   - Row 3 is replaced with a brand-new contract "794/DR KESH" (CIN KS10293,
#     Karim benzima, periodicity "annuelle") with its own monetary figures.
#   - Rows 4-7 (the old "108/ANSYSFYSN01/AV1" contracts) are removed entirely.
#   - The totals row moves from row 8 up to row 4 and its totals are updated
#     to reflect the two remaining contracts.

# Remove the four obsolete data rows (old rows 4-7). Remaining rows below
# shift up, so the old totals row (row 8) becomes row 4.
$ws.Rows("4:7").Delete()

# --- Row 2: contract 002/DR002 ---
$ws.Range("A2").Value = "002/DR002"
$ws.Range("H2").Value = 100000
$ws.Range("I2").Value = 0
$ws.Range("J2").Value = 7500
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 50000
$ws.Range("M2").Value = 92500

# --- Row 3: contract 794/DR KESH (replaces old "Karami abdelilah" row) ---
$ws.Range("A3").Value = "794/DR KESH"
$ws.Range("C3").Value = "KS10293"
$ws.Range("D3").Value = "Karim benzima"
$ws.Range("F3").Value = "annuelle"
$ws.Range("H3").Value = 100000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 10000
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 90000

# --- Row 4: totals row (formerly row 8), updated sums ---
$ws.Range("H4").Value = 200000
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 17500
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 50000
$ws.Range("M4").Value = 182500
